$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B23: prepend an intro line to the existing "rights" text
$ws.Range("B23").Value = "สิทธิของประชาชนผู้เสียภาษีมีดังต่อไปนี้`n1.การผ่อนชำระภาษี`n2.การยื่นอุทธรณ์คัดค้านการประเมินภาษี `n3.ขอทุเลาการชำระภาษีอากรโดยจัดให้มีหลักประกันการชำระหนี้ภาษีอากรค้าง `n4.ของดหรือลดเบี้ยปรับและเงินเพิ่มภาษีอากร `n5.ขอคัดเอกสารหรือขอสำเนาเอกสาร "

# Fill in row 31 with a new tag/response pair (fallback "unknown message" response)
$ws.Range("A31").Value = "unknow-message"
$ws.Range("B31").Value = "ขออภัยค่ะ ระบบอัติโนมัติยังไม่เข้าใจคำถามของคุณ หากต้องการสอบถามข้อมูลด้านใดเกี่ยวกับภาษี สามารถถามคำถามดังต่อไปนี้ได้เลยค่ะ`n1.กฎหมายภาษืคืออะไร`n2.ภาษีมีกี่ประเภท`n3.สิทธิของประชาชนผู้เสียภาษีมีอะไรบ้าง`n4.หน่วยงานของรัฐสามารถจัดเก็บภาษีอากรประเภทไหนได้บ้าง`n5.หากไม่เสียภาษีตามที่กฎหมายกำหนด จะส่งผลอย่างไร`n6.ทำไมหน่วยงานรัฐถึงต้องจัดเก็บภาษี?`n7.ภาษีบุคคลธรรมดาคำนวณจากอะไรและคำนวนอย่างไร`n8.กฎหมายกำหนดให้บุคคลต้องทำการยื่นเสียภาษีในช่วงเดือนใด`n9.ถ้าต้องการที่จะเสียภาษี สามารถยื่นเสียภาษีได้ที่ไหนบ้าง `n10.บุคคลธรรมดาต้องมีเงินเดือนเท่าไร ถึงต้องยื่นภาษี "

# Match the author's final view state: scrolled so row 28 is at top, with B31 selected
$ws.Range("B31").Select()
$ws.Application.ActiveWindow.ScrollRow = 28
